$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Field Explanation": fix/clarify a handful of field-help strings.
# ---------------------------------------------------------------------------
$wsField = $wb.Worksheets.Item("Field Explanation")

# Flags field help text now documents StaffPayed instead of ShowBlip.
$wsField.Range("C9").Value = "Used as 'bulk storage' for the boolean variables [Ownable][Owned][ContextMission][StaffPayed]"

# Example flags value updated to match the real data (1011).
$wsField.Range("C11").Value = "Eg.: 1011"

# Typo fix: Unsed -> Unused.
$wsField.Range("C22").Value = "Unused when flag is already true"

# Typo fix: met -> meet.
$wsField.Range("C34").Value = "If player can't meet StaffSal * Staff, no income is calculated"

# Move the selection to C15, matching the author's last cursor position.
$wsField.Range("C15").Select()

# ---------------------------------------------------------------------------
# Sheet "Properties Table": populate income related columns for Grotti and
# fix the Flags value for Perseus.
# ---------------------------------------------------------------------------
$wsProps = $wb.Worksheets.Item("Properties Table")
$wsProps.Activate()

# Grotti (row 2): IncomeMin, IncomeMax, StaffSal, StaffPay now hold real values.
$wsProps.Range("G2").Value = 100
$wsProps.Range("H2").Value = 1500000
$wsProps.Range("I2").Value = 6500
$wsProps.Range("J2").Value = 1500

# Perseus (row 3): Flags corrected from 1002 to 1001.
$wsProps.Range("C3").Value = 1001

# Select the generated SQL column, matching the author's last selection.
$wsProps.Range("M2:M11").Select()
